$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# --- Row 78: gains a Status ("DONE") -----------------------------------
# (set the value BEFORE hiding the row - doing it after provokes a
#  spurious row-height recalculation in this host)
$ws.Cells.Item(78, 4).Value = "DONE"

# --- New data rows 79-81 -------------------------------------------------
# Row 79 - "display name should have * from model", status DONE
$ws.Cells.Item(79, 1).Value = 79
$ws.Cells.Item(79, 4).Value = "DONE"
$ws.Cells.Item(79, 6).Value = "display name should have * from model"
$ws.Rows.Item(79).RowHeight = 29

# Row 80 - "client side logging" (Title column first so shared-string
# order matches the authored workbook, Description filled in afterwards)
$ws.Cells.Item(80, 1).Value = 80
$ws.Cells.Item(80, 6).Value = "client side logging"

# Row 81 - "server side logging" / "use node debug"
$ws.Cells.Item(81, 1).Value = 81
$ws.Cells.Item(81, 6).Value = "server side logging"
$ws.Cells.Item(81, 9).Value = "use node debug"

# Row 80 Description - "use ngx-logger" (written last on purpose, see above)
$ws.Cells.Item(80, 9).Value = "use ngx-logger"

# --- Hide rows 25 and 78 (now resolved/done issues) ----------------------
$ws.Rows.Item(25).Hidden = $true
$ws.Rows.Item(78).Hidden = $true

# --- AutoFilter: grow the filtered range to the new data, drop the
#     "AWAIT TEST" criterion (keep "PENDING" + blanks) --------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:I78").AutoFilter(4, @("PENDING", ""), 7)

# --- Keep the hidden _FilterDatabase defined name in sync -----------------
$wb.Names.Item("Issues!_FilterDatabase").RefersTo = "=Issues!`$A`$1:`$I`$78"

# --- Restore the active selection to match the authored workbook ---------
$ws.Range("E86").Select()

Write-Host "done"
